# Update cryptocurrency price/volume snapshot values pulled from coinranking.com
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2: Bitcoin
$ws.Range('D2').Value = '75.641.72'
$ws.Range('E2').Value = '  +8.77%  '
# Row 3: Ethereum
$ws.Range('D3').Value = '2.676.69'
$ws.Range('E3').Value = '  +10.02%  '
# Row 4: TetherUSD
$ws.Range('E4').Value = '  +0.15%  '
# Row 5: Solana
$ws.Range('D5').Value = '''187.68'
$ws.Range('E5').Value = '  +13.01%  '
# Row 6: BNB
$ws.Range('D6').Value = '''587.40'
$ws.Range('E6').Value = '  +4.16%  '
# Row 7: USDC
$ws.Range('E7').Value = '  -0.03%  '
# Row 8: XRP
$ws.Range('D8').Value = '''0.539'
$ws.Range('E8').Value = '  +5.23%  '
# Row 9: Dogecoin
$ws.Range('D9').Value = '''0.195'
$ws.Range('E9').Value = '  +14.95%  '
# Row 10: LidoStakedEther
$ws.Range('D10').Value = '2.675.66'
$ws.Range('E10').Value = '  +9.94%  '
# Row 12: Cardano
$ws.Range('E12').Value = '  +7.37%  '
# Row 13: Toncoin
$ws.Range('D13').Value = '''4.72'
$ws.Range('E13').Value = '  +0.82%  '
# Row 14: WrappedBTC
$ws.Range('D14').Value = '75.422.51'
$ws.Range('E14').Value = '  +8.62%  '
# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '3.170.01'
$ws.Range('E15').Value = '  +10.14%  '
# Row 16: ShibaInu
$ws.Range('D16').Value = '''0.0000188'
$ws.Range('E16').Value = '  +5.74%  '
# Row 17: Avalanche
$ws.Range('D17').Value = '''26.53'
$ws.Range('E17').Value = '  +10.77%  '
# Row 18: WrappedEther
$ws.Range('D18').Value = '2.678.47'
$ws.Range('E18').Value = '  +10.34%  '
# Row 19: Uniswap
$ws.Range('D19').Value = '''9.26'
$ws.Range('E19').Value = '  +29.72%  '
# Row 20: Chainlink
$ws.Range('D20').Value = '''11.96'
$ws.Range('E20').Value = '  +10.84%  '
# Row 21: BitcoinCash
$ws.Range('D21').Value = '''372.61'
$ws.Range('E21').Value = '  +9.14%  '
# Row 22: SuiNetwork
$ws.Range('D22').Value = '''2.28'
$ws.Range('E22').Value = '  +15.72%  '
# Row 23: Polkadot
$ws.Range('D23').Value = '''4.09'
$ws.Range('E23').Value = '  +5.19%  '
# Row 24: LEO
$ws.Range('E24').Value = '  +4.46%  '
# Row 26: Litecoin
$ws.Range('E26').Value = '  +6.44%  '
# Row 27: NEARProtocol
$ws.Range('E27').Value = '  +9.56%  '
# Row 28: Aptos
$ws.Range('D28').Value = '''9.37'
$ws.Range('E28').Value = '  +9.88%  '
# Row 29: WrappedeETH
$ws.Range('D29').Value = '2.820.22'
$ws.Range('E29').Value = '  +10.45%  '
# Row 30: Binance-PegBSC-USD
$ws.Range('E30').Value = '  +0.59%  '
# Row 31: PEPE
$ws.Range('D31').Value = '0.0₃0945'
$ws.Range('E31').Value = '  +11.22%  '
# Row 32: Fetch.AI
$ws.Range('D32').Value = '''1.41'
$ws.Range('E32').Value = '  +14.18%  '
# Row 33: Bittensor
$ws.Range('D33').Value = '''516.52'
$ws.Range('E33').Value = '  +14.28%  '
# Row 34: InternetComputer(DFINITY)
$ws.Range('D34').Value = '''7.74'
$ws.Range('E34').Value = '  +4.76%  '
# Row 35: PancakeSwap
$ws.Range('D35').Value = '''1.75'
$ws.Range('E35').Value = '  +8.82%  '
# Row 36: FirstDigitalUSD
$ws.Range('D36').Value = '''1.00'
$ws.Range('E36').Value = '  +0.10%  '
# Row 37: Monero
$ws.Range('D37').Value = '''163.48'
$ws.Range('E37').Value = '  +3.82%  '
# Row 38: Kaspa
$ws.Range('E38').Value = '  +7.81%  '
# Row 39: EthereumClassic
$ws.Range('D39').Value = '''19.22'
$ws.Range('E39').Value = '  +5.56%  '
# Row 40: WhiteBITCoin
$ws.Range('D40').Value = '''19.39'
$ws.Range('E40').Value = '  +1.45%  '
# Row 42: RenderToken
$ws.Range('D42').Value = '''5.01'
$ws.Range('E42').Value = '  +14.04%  '
# Row 43: Aave
$ws.Range('D43').Value = '''168.92'
$ws.Range('E43').Value = '  +25.19%  '
# Row 44: Stacks
$ws.Range('D44').Value = '''1.70'
$ws.Range('E44').Value = '  +12.05%  '
# Row 45: PolygonEcosystemToken
$ws.Range('D45').Value = '''0.331'
$ws.Range('E45').Value = '  +8.78%  '
# Row 46: ImmutableX
$ws.Range('D46').Value = '''1.19'
$ws.Range('E46').Value = '  +9.69%  '
# Row 47: dogwifhat
$ws.Range('D47').Value = '''2.37'
$ws.Range('E47').Value = '  +13.19%  '
# Row 48: OKB
$ws.Range('D48').Value = '''39.27'
$ws.Range('E48').Value = '  +3.57%  '
# Row 49: Cronos
$ws.Range('D49').Value = '''0.0843'
$ws.Range('E49').Value = '  +16.45%  '
# Row 50: Filecoin
$ws.Range('E50').Value = '  +7.45%  '
# Row 51: ARBITRUM
$ws.Range('D51').Value = '''0.535'
$ws.Range('E51').Value = '  +9.53%  '
